$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.070.07'
$ws.Range("E2").Value = '  +0.98%  '
$ws.Range("D3").Value = '1.741.67'
$ws.Range("E3").Value = '  +0.05%  '
$ws.Range("E4").Value = '  +0.01%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '241.79'
$c.ClearFormats()
$ws.Range("E5").Value = '  +4.66%  '
$ws.Range("E6").Value = '  +0.07%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.5280'
$c.ClearFormats()
$ws.Range("E7").Value = '  +2.26%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.2783'
$c.ClearFormats()
$ws.Range("E8").Value = '  -0.83%  '
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.06177'
$c.ClearFormats()
$ws.Range("E9").Value = '  +1.19%  '
$ws.Range("B10").Value = 'WrappedEther'
$ws.Range("C10").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D10").Value = '1.740.45'
$ws.Range("E10").Value = '  -0.14%  '
$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.07193'
$c.ClearFormats()
$ws.Range("E11").Value = '  +2.15%  '
$ws.Range("B12").Value = 'Solana'
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '15.26'
$c.ClearFormats()
$ws.Range("E12").Value = '  -0.15%  '
$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.6523'
$c.ClearFormats()
$ws.Range("E13").Value = '  +1.69%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '4.630'
$c.ClearFormats()
$ws.Range("E14").Value = '  +2.54%  '
$ws.Range("B15").Value = 'Litecoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '77.85'
$c.ClearFormats()
$ws.Range("E15").Value = '  +1.06%  '
$ws.Range("B16").Value = 'Dai'
$ws.Range("C16").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.ClearFormats()
$ws.Range("E16").Value = '  +0.05%  '
$ws.Range("B17").Value = 'BinanceUSD'
$ws.Range("C17").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.ClearFormats()
$ws.Range("E17").Value = '  +0.06%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '26.076.35'
$ws.Range("E18").Value = '  +0.98%  '
$ws.Range("B19").Value = 'Avalanche'
$ws.Range("C19").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '11.87'
$c.ClearFormats()
$ws.Range("E19").Value = '  +3.79%  '
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '0.000006769'
$c.ClearFormats()
$ws.Range("E20").Value = '  +2.87%  '
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '1.964.22'
$ws.Range("E21").Value = '  -0.61%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '4.370'
$c.ClearFormats()
$ws.Range("E22").Value = '  +5.88%  '
$ws.Range("B23").Value = 'Cosmos'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '8.693'
$c.ClearFormats()
$ws.Range("E23").Value = '  +0.70%  '
$ws.Range("B24").Value = 'Chainlink'
$ws.Range("C24").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '5.256'
$c.ClearFormats()
$ws.Range("E24").Value = '  +2.09%  '
$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '140.02'
$c.ClearFormats()
$ws.Range("E25").Value = '  +0.14%  '
$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '1.515'
$c.ClearFormats()
$ws.Range("E26").Value = '  +0.67%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '15.29'
$c.ClearFormats()
$ws.Range("E27").Value = '  +1.38%  '
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '1.771'
$c.ClearFormats()
$ws.Range("E28").Value = '  -2.51%  '
$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '106.00'
$c.ClearFormats()
$ws.Range("E29").Value = '  +3.40%  '
$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '0.08451'
$c.ClearFormats()
$ws.Range("E30").Value = '  +2.21%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '3.863'
$c.ClearFormats()
$ws.Range("E31").Value = '  +5.43%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '3.664'
$c.ClearFormats()
$ws.Range("E32").Value = '  +7.02%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.04612'
$c.ClearFormats()
$ws.Range("E33").Value = '  +2.71%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '2.658'
$c.ClearFormats()
$ws.Range("E34").Value = '  +1.71%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.9969'
$c.ClearFormats()
$ws.Range("E35").Value = '  +1.44%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.6279'
$c.ClearFormats()
$ws.Range("E36").Value = '  +2.06%  '
$ws.Range("B37").Value = 'MXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '2.700'
$c.ClearFormats()
$ws.Range("E37").Value = '  +1.81%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.01612'
$c.ClearFormats()
$ws.Range("E38").Value = '  +1.68%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '1.933'
$c.ClearFormats()
$ws.Range("E39").Value = '  -0.02%  '
$ws.Range("B40").Value = 'PaxDollar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.ClearFormats()
$ws.Range("E40").Value = '  +0.13%  '
$ws.Range("B41").Value = 'Quant'
$ws.Range("C41").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '99.26'
$c.ClearFormats()
$ws.Range("E41").Value = '  -1.18%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.3906'
$c.ClearFormats()
$ws.Range("E42").Value = '  +1.79%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.7547'
$c.ClearFormats()
$ws.Range("E43").Value = '  +4.16%  '
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '4.981'
$c.ClearFormats()
$ws.Range("E44").Value = '  +0.37%  '
$ws.Range("B45").Value = 'Algorand'
$ws.Range("C45").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.1149'
$c.ClearFormats()
$ws.Range("E45").Value = '  +2.30%  '
$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '6.281'
$c.ClearFormats()
$ws.Range("E46").Value = '  +0.12%  '
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.05327'
$c.ClearFormats()
$ws.Range("E47").Value = '  -1.55%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '54.85'
$c.ClearFormats()
$ws.Range("E48").Value = '  +3.14%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '30.82'
$c.ClearFormats()
$ws.Range("E49").Value = '  +3.34%  '
$ws.Range("B50").Value = 'Decentraland'
$ws.Range("C50").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.3473'
$c.ClearFormats()
$ws.Range("E50").Value = '  +2.54%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '7.592'
$c.ClearFormats()
$ws.Range("E51").Value = '  -1.15%  '
